# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) figures across the four worksheets.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) sheet ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7990
$ws1.Range("F4").Value = 91
$ws1.Range("F5").Value = 22774
$ws1.Range("F6").Value = 47
$ws1.Range("F8").Value = 664
$ws1.Range("F10").Value = 142
$ws1.Range("F12").Value = 791
$ws1.Range("F13").Value = 46
$ws1.Range("F14").Value = 594
$ws1.Range("F15").Value = 363
$ws1.Range("F17").Value = 311
$ws1.Range("F19").Value = 403
$ws1.Range("F20").Value = 414
$ws1.Range("F21").Value = 1108
$ws1.Range("F23").Value = 669
$ws1.Range("F24").Value = 2277
$ws1.Range("F25").Value = 788
$ws1.Range("F27").Value = 1057
$ws1.Range("F28").Value = 39
$ws1.Range("F30").Value = 1063

# ---- 演出 (Performance) sheet ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 335

# ---- 本地生活 (Local Life) sheet ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 504

# ---- 全部类型 (All Types) sheet ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 504
$ws4.Range("F3").Value = 7990
$ws4.Range("F5").Value = 91
$ws4.Range("F7").Value = 22820
$ws4.Range("F8").Value = 47
$ws4.Range("F10").Value = 664
$ws4.Range("F13").Value = 142
$ws4.Range("F15").Value = 335
$ws4.Range("F18").Value = 791
$ws4.Range("F19").Value = 46
$ws4.Range("F20").Value = 594
$ws4.Range("F21").Value = 363
$ws4.Range("F27").Value = 311
$ws4.Range("F29").Value = 403
$ws4.Range("F30").Value = 414
$ws4.Range("F31").Value = 1108
$ws4.Range("F33").Value = 669
$ws4.Range("F34").Value = 2277
$ws4.Range("F35").Value = 788
$ws4.Range("F37").Value = 1059
$ws4.Range("F38").Value = 39
$ws4.Range("F41").Value = 1063
